$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 43.314413611233547
$ws.Range("B3").Value = 25.848364715413481
$ws.Range("B4").Value = 33.712925354337223
$ws.Range("B5").Value = 19.619408152239028
$ws.Range("B6").Value = 29.692474850554529
$ws.Range("B7").Value = 19.890954420405809
$ws.Range("B8").Value = 12.636187994095749
$ws.Range("B9").Value = 4.9570439882781852

$ws.Range("B2:B9").Select()
